$p = $ppt.ActivePresentation

# ---------------------------------------------------------------------
# 1) Update the "datetimeFigureOut" date placeholder (02/06/2018 ->
#    08/06/2018) everywhere it appears: the Slide Master and every
#    Slide Layout ("CustomLayout").
# ---------------------------------------------------------------------
function Update-DateFieldShapes($shapes) {
    for ($i = 1; $i -le $shapes.Count; $i++) {
        $sh = $shapes.Item($i)
        if ($sh.HasTextFrame -eq -1) {
            if ($sh.TextFrame.TextRange.Text -eq "02/06/2018") {
                $sh.TextFrame.TextRange.Text = "08/06/2018"
            }
        }
    }
}

$design = $p.Designs.Item(1)
$master = $design.SlideMaster

Update-DateFieldShapes $master.Shapes

for ($li = 1; $li -le $master.CustomLayouts.Count; $li++) {
    $layout = $master.CustomLayouts.Item($li)
    Update-DateFieldShapes $layout.Shapes
}

# ---------------------------------------------------------------------
# 2) Per-slide fixes on the bottom navigation bar shapes (present,
#    with identical names, on every slide of the deck).
#    - "Retângulo de cantos arredondados 4" -> the URL textbox; its
#      "http://" + "www...." runs got split by a stray edit and need
#      to be merged back into one run.
#    - "Retângulo de cantos arredondados 6" -> the footer bar; it was
#      resized/repositioned (made thinner, moved down).
#    - "Retângulo 9" -> nav button; stale "Vendas" label renamed to
#      "Reservas" to match the other slides.
# ---------------------------------------------------------------------
for ($si = 1; $si -le $p.Slides.Count; $si++) {
    $s = $p.Slides.Item($si)

    # -- merge the split "http://" + "www...." runs back into a single
    #    run (keeps the first run's formatting)
    $urlShape = $s.Shapes.Item("Retângulo de cantos arredondados 4")
    $tr = $urlShape.TextFrame.TextRange
    $full = $tr.Text
    if ($full.StartsWith("http://")) {
        $len = $tr.Length
        $whole = $tr.Characters(1, $len)
        $whole.Text = $full
    }

    # -- fix the footer rounded-rectangle bar position/size
    $barShape = $s.Shapes.Item("Retângulo de cantos arredondados 6")
    $barShape.Top = 519.4765354330709
    $barShape.Height = 20.525354330708662

    # -- rename the stale "Vendas" nav button to "Reservas"
    $navShape = $s.Shapes.Item("Retângulo 9")
    $navTr = $navShape.TextFrame.TextRange
    if ($navTr.Text -eq "Vendas") {
        $navLen = $navTr.Length
        $navWhole = $navTr.Characters(1, $navLen)
        $navWhole.Text = "Reservas"
    }
}
